$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-column style (s="2": bold border + center/top align + YYYY-MM-DD HH:MM:SS numfmt)
# from the last existing data row (385) down through the new rows (386-464) so no new
# style entries get created in styles.xml.
$ws.Range("A385").Copy() | Out-Null
$ws.Range("A386:A464").PasteSpecial(-4122) | Out-Null

$data = @(
  @(44460, 0, 8, 86.16975441619991),
  @(44461, 0, 8, 86.16975441619991),
  @(44462, 1, 5, 53.85609651012494),
  @(44463, 0, 5, 53.85609651012494),
  @(44464, 2, 7, 75.39853511417492),
  @(44465, 0, 3, 32.31365790607497),
  @(44466, 0, 3, 32.31365790607497),
  @(44467, 0, 3, 32.31365790607497),
  @(44468, 0, 3, 32.31365790607497),
  @(44469, 0, 2, 21.54243860404998),
  @(44470, 0, 2, 21.54243860404998),
  @(44471, 0, 0, 0),
  @(44472, 1, 1, 10.77121930202499),
  @(44473, 0, 1, 10.77121930202499),
  @(44474, 0, 1, 10.77121930202499),
  @(44475, 2, 3, 32.31365790607497),
  @(44476, 1, 4, 43.08487720809995),
  @(44477, 1, 5, 53.85609651012494),
  @(44478, 1, 6, 64.62731581214993),
  @(44479, 0, 5, 53.85609651012494),
  @(44480, 1, 6, 64.62731581214993),
  @(44481, 0, 6, 64.62731581214993),
  @(44482, 0, 4, 43.08487720809995),
  @(44483, 0, 3, 32.31365790607497),
  @(44484, 0, 2, 21.54243860404998),
  @(44485, 0, 1, 10.77121930202499),
  @(44486, 0, 1, 10.77121930202499),
  @(44487, 0, 0, 0),
  @(44488, 0, 0, 0),
  @(44489, 1, 1, 10.77121930202499),
  @(44490, 0, 1, 10.77121930202499),
  @(44491, 1, 2, 21.54243860404998),
  @(44492, 0, 2, 21.54243860404998),
  @(44493, 0, 2, 21.54243860404998),
  @(44494, 3, 5, 53.85609651012494),
  @(44495, 0, 5, 53.85609651012494),
  @(44496, 0, 4, 43.08487720809995),
  @(44497, 3, 7, 75.39853511417492),
  @(44498, 1, 7, 75.39853511417492),
  @(44499, 2, 9, 96.9409737182249),
  @(44500, 10, 19, 204.6531667384748),
  @(44501, 2, 18, 193.8819474364498),
  @(44502, 11, 29, 312.3653597587247),
  @(44503, 0, 29, 312.3653597587247),
  @(44504, 0, 26, 280.0517018526497),
  @(44505, 6, 31, 333.9077983627747),
  @(44506, 4, 33, 355.4502369668246),
  @(44507, 1, 24, 258.5092632485997),
  @(44508, 6, 28, 301.5941404566997),
  @(44509, 0, 17, 183.1107281344248),
  @(44510, 0, 17, 183.1107281344248),
  @(44511, 2, 19, 204.6531667384748),
  @(44512, 8, 21, 226.1956053425248),
  @(44513, 2, 19, 204.6531667384748),
  @(44514, 3, 21, 226.1956053425248),
  @(44515, 2, 17, 183.1107281344248),
  @(44516, 12, 29, 312.3653597587247),
  @(44517, 0, 29, 312.3653597587247),
  @(44518, 1, 28, 301.5941404566997),
  @(44519, 2, 22, 236.9668246445498),
  @(44520, 3, 23, 247.7380439465747),
  @(44521, 0, 20, 215.4243860404998),
  @(44522, 2, 20, 215.4243860404998),
  @(44523, 2, 10, 107.7121930202499),
  @(44524, 2, 12, 129.2546316242999),
  @(44525, 0, 11, 118.4834123222749),
  @(44526, 3, 12, 129.2546316242999),
  @(44527, 1, 10, 107.7121930202499),
  @(44528, 2, 12, 129.2546316242999),
  @(44529, 0, 10, 107.7121930202499),
  @(44530, 1, 9, 96.9409737182249),
  @(44531, 0, 7, 75.39853511417492),
  @(44532, 1, 8, 86.16975441619991),
  @(44533, 0, 5, 53.85609651012494),
  @(44534, 2, 6, 64.62731581214993),
  @(44535, 1, 5, 53.85609651012494),
  @(44536, 1, 6, 64.62731581214993),
  @(44537, 0, 5, 53.85609651012494),
  @(44538, 1, 6, 64.62731581214993),
)

$row = 386
foreach ($rec in $data) {
  $ws.Cells.Item($row, 1).Value = $rec[0]
  $ws.Cells.Item($row, 2).Value = $rec[1]
  $ws.Cells.Item($row, 3).Value = $rec[2]
  $ws.Cells.Item($row, 4).Value = $rec[3]
  $row = $row + 1
}
